# Update countries & provincias Spain
# Applies the 3-April-2020 20:20 data refresh: reorders several countries
# in the list (new countries inserted alphabetically/positionally) and
# refreshes the numeric COVID stats for every row, plus the title timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 20:20"

$data = @(
  @("Estados Unidos",265506,20629,11983,246737,5781,716,6786),
  @("Italia",119827,4585,19758,85388,4068,766,14681),
  @("España",117710,5645,30513,76262,6416,587,10935),
  @("Alemania",89838,5044,24575,64033,3936,123,1230),
  @("China",81620,31,76571,1727,379,4,3322),
  @("Francia",64338,5233,14008,43823,6662,1120,6507),
  @("Iran",53183,2715,17935,31954,4035,134,3294),
  @("Reino Unido",38168,4450,135,34428,163,684,3605),
  @("Turquia",20921,2786,484,20012,1251,69,425),
  @("Suiza",19606,779,4846,14169,348,55,591),
  @("Belgica",16770,1422,2872,12755,1205,132,1143),
  @("Paises Bajos",15723,1026,250,13986,1182,148,1487),
  @("Canada",11747,464,1979,9595,120,0,173),
  @("Austria",11506,377,2022,9316,245,10,168),
  @("Corea del Sur",10062,86,6021,3867,55,5,174),
  @("Portugal",9886,852,68,9572,245,37,246),
  @("Brasil",8229,185,127,7759,296,19,343),
  @("Israel",7428,571,338,7050,113,4,40),
  @("Suecia",6131,563,205,5568,469,50,358),
  @("Australia",5350,36,585,4737,85,3,28),
  @("Noruega",5296,149,32,5207,96,7,57),
  @("Irlanda",4273,424,5,4148,109,22,120),
  @("Rusia",4149,601,281,3834,8,4,34),
  @("Chequia",4091,233,72,3966,77,9,53),
  @("Dinamarca",3757,371,1193,2425,153,16,139),
  @("Chile",3737,333,427,3288,31,4,22),
  @("Ecuador",3368,205,65,3158,100,25,145),
  @("Malasia",3333,217,827,2453,108,3,53),
  @("Polonia",3266,320,56,3145,50,8,65),
  @("Rumania",3183,445,283,2767,83,18,133),
  @("Filipinas",3018,385,52,2830,1,29,136),
  @("Pakistan",2637,216,126,2471,10,6,40),
  @("Japon",2617,0,514,2040,60,0,63),
  @("Luxemburgo",2612,125,500,2081,33,1,31),
  @("India",2567,24,192,2303,0,0,72),
  @("Arabia Saudita",2039,154,351,1663,41,4,25),
  @("Indonesia",1986,196,134,1671,0,11,181),
  @("Tailandia",1978,103,581,1378,23,4,19),
  @("Finlandia",1615,97,300,1295,72,1,20),
  @("Grecia",1613,69,61,1493,92,6,59),
  @("Peru",1595,181,537,1003,51,0,55),
  @("Mexico",1510,132,633,827,1,13,50),
  @("Sudafrica",1505,43,95,1401,7,4,9),
  @("Republica Dominicana",1488,108,16,1404,147,8,68),
  @("Serbia",1476,305,42,1395,81,8,39),
  @("Panama",1475,0,10,1428,50,0,37),
  @("Islandia",1364,45,309,1051,12,0,4),
  @("Argentina",1265,0,266,960,0,3,39),
  @("Emiratos Arabes Unidos",1264,240,108,1147,2,1,9),
  @("Argelia",1171,185,62,1004,0,19,105),
  @("Colombia",1161,0,55,1087,50,0,19),
  @("Singapur",1114,65,282,827,24,1,5),
  @("Croacia",1079,68,92,979,39,1,8),
  @("Catar",1075,126,93,979,37,0,3),
  @("Estonia",961,103,48,901,16,1,12),
  @("Ucrania",942,45,19,900,16,1,23),
  @("Eslovenia",934,37,70,844,31,3,20),
  @("Nueva Zelanda",868,71,103,764,2,0,1),
  @("Egipto",865,0,201,606,0,0,58),
  @("Hong Kong",845,43,173,668,8,0,4),
  @("Irak",820,48,226,540,0,0,54),
  @("Marruecos",761,53,56,658,1,3,47),
  @("Armenia",736,73,43,686,30,0,7),
  @("Crucero",712,0,619,82,10,0,11),
  @("Lituania",696,47,7,680,11,0,9),
  @("Barein",672,29,382,286,3,0,4),
  @("Hungria",623,38,43,554,17,5,26),
  @("Moldavia",591,86,26,557,65,2,8),
  @("Bosnia y Herzegovina",575,42,27,531,4,1,17),
  @("Libano",508,14,50,441,26,1,17),
  @("Tunez",495,40,5,472,10,4,18),
  @("Letonia",493,35,1,491,3,1,1),
  @("Bulgaria",485,28,30,441,18,4,14),
  @("Kazajistan",460,25,29,425,6,3,6),
  @("Eslovaquia",450,24,10,439,3,0,1),
  @("Azerbaiyan",443,43,32,406,7,0,5),
  @("Principado de Andorra",439,11,16,407,12,1,16),
  @("Republica de Macedonia",430,46,17,402,8,0,11),
  @("Kuwait",417,75,82,335,16,0,0),
  @("Costa Rica",396,0,6,388,11,0,2),
  @("Republica de Chipre",396,40,28,357,11,1,11),
  @("Uruguay",369,19,68,297,13,0,4),
  @("Bielorrusia",351,47,53,294,11,0,4),
  @("Taiwan",348,9,50,293,0,0,5),
  @("Reunion",321,13,40,281,3,0,0),
  @("Jordania",310,11,58,247,5,0,5),
  @("Camerun",306,0,10,289,0,0,7),
  @("Albania",304,27,89,198,7,1,17),
  @("Burkina Faso",302,14,50,236,0,0,16),
  @("Afganistan",281,8,10,265,0,0,6),
  @("Cuba",269,36,15,248,8,0,6),
  @("Oman",252,21,57,194,3,0,1),
  @("San Marino",245,0,21,194,15,0,30),
  @("Vietnam",237,4,85,152,3,0,0),
  @("Honduras",222,3,3,204,10,1,15),
  @("Uzbekistan",221,16,25,194,8,0,2),
  @("Senegal",207,12,66,140,1,0,1),
  @("Ghana",205,1,31,169,2,0,5),
  @("Malta",202,6,2,200,2,0,0),
  @("Costa de Marfil",194,0,15,178,0,0,1),
  @("Estado de Palestina",193,32,21,171,0,0,1),
  @("Nigeria",190,6,20,168,0,0,2),
  @("Mauricio",186,17,0,179,1,0,7),
  @("Islas Feroe",179,2,91,88,1,0,0),
  @("Montenegro",174,30,1,171,4,0,2),
  @("Sri Lanka",159,8,24,131,5,0,4),
  @("Georgia",148,14,27,121,6,0,0),
  @("Venezuela",146,0,43,98,6,0,5),
  @("Martinica",138,0,27,108,19,0,3),
  @("Consejo Danes para los Refugiados",134,0,3,118,0,0,13),
  @("Brunei",134,1,65,68,3,0,1),
  @("Bolivia",132,9,1,122,3,1,9),
  @("Kirguistan",130,14,6,123,5,1,1),
  @("Guadalupe",128,0,24,98,14,0,6),
  @("Kenia",122,12,4,114,2,1,4),
  @("Mayotte",116,0,10,105,3,0,1),
  @("Isla de Man",114,19,0,113,0,0,1),
  @("Camboya",114,4,35,79,1,0,0),
  @("Niger",98,0,0,93,0,0,5),
  @("Trinidad yTobago",97,3,1,90,0,1,6),
  @("Gibraltar",95,7,46,49,0,0,0),
  @("Paraguay",92,15,6,83,2,0,3),
  @("Ruanda",84,0,0,84,0,0,0),
  @("Liechtenstein",75,0,0,75,0,0,0),
  @("Guinea",73,21,2,71,0,0,0),
  @("Madagascar",65,6,0,65,6,0,0),
  @("Monaco",64,4,3,60,2,0,1),
  @("Aruba",62,2,1,61,0,0,0),
  @("Banglades",61,5,26,29,1,0,6),
  @("Guayana Francesa",57,6,22,35,1,0,0),
  @("Barbados",51,5,0,51,0,0,0),
  @("Guatemala",50,3,12,37,1,0,1),
  @("Republica de Yibuti",49,9,8,41,0,0,0),
  @("Jamaica",47,0,2,42,0,0,3),
  @("El Salvador",46,5,0,44,4,0,2),
  @("Uganda",45,0,0,45,0,0,0),
  @("Macao",42,1,10,32,0,0,0),
  @("Togo",40,1,17,20,0,1,3),
  @("Zambia",39,0,2,36,0,0,1),
  @("Puerto Rico",39,0,1,36,0,0,2),
  @("Polinesia Francesa",37,0,0,37,1,0,0),
  @("Mali",36,0,0,33,0,0,3),
  @("Etiopia",35,6,3,32,2,0,0),
  @("Bermudas",35,0,11,24,0,0,0),
  @("Guam",32,0,0,31,0,0,1),
  @("Islas Caimanes",28,0,0,27,0,0,1),
  @("Bahamas",24,0,1,22,1,0,1),
  @("San Martin (Parte Holandesa)",23,5,6,15,0,1,2),
  @("Eritrea",22,0,0,22,0,0,0),
  @("San Martin (Parte Francesa)",22,0,2,19,0,0,1),
  @("Congo",22,0,2,18,0,0,2),
  @("Gabon",21,0,1,19,0,0,1),
  @("Birmania",20,0,0,19,0,0,1),
  @("Tanzania",20,0,3,16,0,0,1),
  @("Guyana",19,0,0,15,0,0,4),
  @("Maldivas",19,0,13,6,0,0,0),
  @("Nueva Caledonia",18,0,1,17,0,0,0),
  @("Haiti",18,2,1,17,0,0,0),
  @("Islas Virgenes de los Estados Unidos",17,0,0,17,0,0,0),
  @("Guinea Ecuatorial",16,1,1,15,0,0,0),
  @("Siria",16,0,0,14,0,0,2),
  @("Mongolia",14,0,2,12,0,0,0),
  @("Namibia",14,0,3,11,0,0,0),
  @("Santa Lucia",13,0,1,12,0,0,0),
  @("Benin",13,0,1,12,0,0,0),
  @("Dominica",12,0,0,12,0,0,0),
  @("Libia",11,0,0,10,0,0,1),
  @("Curazao",11,0,3,7,0,0,1),
  @("Mozambique",10,0,0,10,0,0,0),
  @("Seychelles",10,0,0,10,0,0,0),
  @("Granada",10,0,0,10,0,0,0),
  @("Laos",10,0,0,10,0,0,0),
  @("Surinam",10,0,0,9,0,1,1),
  @("Groenlandia",10,0,3,7,0,0,0),
  @("Sudan",10,2,2,6,0,0,2),
  @("Guinea-Bisau",9,0,0,9,0,0,0),
  @("Antigua y Barbuda",9,0,0,9,0,0,0),
  @("Suazilandia",9,0,0,9,0,0,0),
  @("San Cristobal y Nieves",9,0,0,9,0,0,0),
  @("Zimbabue",9,0,0,8,0,0,1),
  @("Montserrat",9,0,0,7,0,0,2),
  @("Republica del Chad",8,0,0,8,0,0,0),
  @("Republica de Africa Central",8,5,0,8,0,0,0),
  @("Angola",8,0,1,5,0,0,2),
  @("Liberia",7,1,0,7,0,0,0),
  @("Fiyi",7,0,0,7,0,0,0),
  @("Santa Sede",7,0,0,7,0,0,0),
  @("Nepal",6,0,1,5,0,0,0),
  @("San Bartolome",6,0,1,5,0,0,0),
  @("Cabo Verde",6,0,0,5,0,0,1),
  @("Mauritania",6,0,2,3,0,0,1),
  @("Islas Turcas y Caicos",5,0,0,5,0,0,0),
  @("Somalia",5,0,1,4,0,0,0),
  @("Nicaragua",5,0,0,4,0,0,1),
  @("Butan",5,0,2,3,0,0,0),
  @("Belice",4,1,0,4,0,0,0),
  @("Botsuana",4,0,0,3,0,0,1),
  @("Gambia",4,0,2,1,0,0,1),
  @("Anguila",3,0,0,3,0,0,0),
  @("Islas Virgenes Britanicas",3,0,0,3,0,0,0),
  @("Malaui",3,0,0,3,0,0,0),
  @("Burundi",3,0,0,3,0,0,0),
  @("San Vicente y las Granadinas",3,1,1,2,0,0,0),
  @("Sierra Leona",2,0,0,2,0,0,0),
  @("Bonaire, San Eustaquio y Saba",2,0,0,2,0,0,0),
  @("Papua Nueva Guinea",1,0,0,1,0,0,0),
  @("Timor Oriental",1,0,0,1,0,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = 4 + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
